# "added new data for Policia" - add a new "2025" sheet with the latest
# victim-by-age/gender breakdown, placed after the "2024" sheet.

$wb = $excel.ActiveWorkbook

$sheet2024 = $wb.Worksheets.Item("2024")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet2024)
$newSheet.Name = "2025"

# Header row
$newSheet.Range("A1").Value = "Grupos de Edad"
$newSheet.Range("B1").Value = "Mujeres"
$newSheet.Range("C1").Value = "Hombres"

# Data rows
$newSheet.Range("A2").Value = "Menores"
$newSheet.Range("B2").Value = 943
$newSheet.Range("C2").Value = 393

$newSheet.Range("A3").Value = "Mayores"
$newSheet.Range("B3").Value = 247
$newSheet.Range("C3").Value = 50

$newSheet.Range("A4").Value = "Desconocida"
$newSheet.Range("B4").Value = 73
$newSheet.Range("C4").Value = 27

# Widen column A so the labels aren't truncated.
$newSheet.Columns.Item(1).ColumnWidth = 15.3

# The previously-active "2024" sheet is left with A1:C4 selected.
$sheet2024.Activate() | Out-Null
$sheet2024.Range("A1:C4").Select() | Out-Null

# New sheet becomes the active tab, with B2:C4 selected.
$newSheet.Activate() | Out-Null
$newSheet.Range("B2:C4").Select() | Out-Null
